# Insert a new "WSTG-INPV-21 / Testing for CSV Injection" row into the
# "Testing Checklist" sheet, directly above the existing blank separator
# row that currently sits at row 95 (between the Input Validation Testing
# section and the Testing for Error Handling section).
#
# Net effect matches the upstream OOXML diff:
#   - dimension grows from A1:G139 to A1:G140
#   - everything from the old row 95 down to row 139 shifts down by one row
#   - the new row 95 carries the WSTG-INPV-21 content
#   - the B4:F139 conditional formatting range grows to B4:F140
#   - the E-column "status" data validation list picks up the new E95 cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing Checklist")

# --- 1. Insert a blank row at 95, shifting rows 95-139 down to 96-140 ---
$ws.Rows.Item(95).Insert()

# --- 2. Give the new row the same look as the other WSTG test-item rows ---
# (Row 94 - WSTG-INPV-20 - uses the standard "test row" formatting that the
# new row should also use.)
$ws.Range("A94:F94").Copy()
$ws.Range("A95:F95").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Rows.Item(95).RowHeight = 181.5

# --- 3. Fill in the new row's content ---
$ws.Range("B95").Value = "WSTG-INPV-21"

$ws.Range("C95").Formula = "=HYPERLINK(""https://owasp.org/www-project-web-security-testing-guide/latest/4-Web_Application_Security_Testing/07-Input_Validation_Testing/21-Testing_for_CSV_Injection"", ""Testing for CSV Injection"")"

$ws.Range("D95").Value = "- Identify CSV/spreadsheet export features that include untrusted input.
- Verify whether attacker-controlled values are interpreted as formulas when the export is opened in common spreadsheet applications.
- Check whether separator/quote injection can move a dangerous prefix to the start of a cell.
- Validate whether mitigations remain effective in Microsoft Excel after saving and re-opening the CSV.
- Assess practical impact based on who opens the export and how it is used."

$ws.Range("E95").Value = "Not Started"

# --- 4. Extend the B4:F139 conditional formatting block down to B4:F140 ---
$fcs = $ws.Range("B100").FormatConditions()
$fcCount = $fcs.Count()
for ($i = 1; $i -le $fcCount; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("B4:F140"))
}

# --- 5. Add the "Not Started,Pass,Issues,N/A" list validation to E95 ---
$ws.Range("E95").Validation.Add(3, 1, 1, "Not Started,Pass,Issues,N/A")
$ws.Range("E95").Validation.ShowInput = $false
$ws.Range("E95").Validation.ShowError = $false
